# Auto-generated Excel COM-interop script to apply the Masamune_Profits data refresh
# (scheduled runner update: currentAveragePrice / LevePrice / LeveProfit columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 521.5
$ws.Range("I4").Value = 383.5
$ws.Range("J4").Value = 797.5
$ws.Range("K4").Value = 383.5
$ws.Range("L4").Value = 797.5
$ws.Range("M4").Value = -269.5
$ws.Range("N4").Value = -1025.5
$ws.Range("H15").Value = 351.92
$ws.Range("I15").Value = 351.92
$ws.Range("K15").Value = 1055.76
$ws.Range("M15").Value = -886.76
$ws.Range("H18").Value = 606.8570999999999
$ws.Range("I18").Value = 462.5
$ws.Range("J18").Value = 799.3333
$ws.Range("K18").Value = 462.5
$ws.Range("L18").Value = 799.3333
$ws.Range("M18").Value = -178.5
$ws.Range("N18").Value = -1367.3333
$ws.Range("H93").Value = 31967.924
$ws.Range("J93").Value = 31967.924
$ws.Range("L93").Value = 31967.924
$ws.Range("N93").Value = -36959.924
$ws.Range("H99").Value = 1987
$ws.Range("I99").Value = 1987
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5961
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -4463
$ws.Range("N99").ClearContents()
$ws.Range("H109").Value = 39675
$ws.Range("J109").Value = 39675
$ws.Range("L109").Value = 39675
$ws.Range("N109").Value = -42449
$ws.Range("H116").Value = 6842.2104
$ws.Range("I116").Value = 3875
$ws.Range("J116").Value = 9000.182000000001
$ws.Range("K116").Value = 3875
$ws.Range("L116").Value = 9000.182000000001
$ws.Range("M116").Value = -433
$ws.Range("N116").Value = -15884.182
$ws.Range("H128").Value = 41167.668
$ws.Range("J128").Value = 41167.668
$ws.Range("L128").Value = 41167.668
$ws.Range("N128").Value = -51127.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2978.0476
$ws.Range("I61").Value = 1771
$ws.Range("K61").Value = 1771
$ws.Range("M61").Value = -1559
$ws.Range("H80").Value = 52997.332
$ws.Range("J80").Value = 52997.332
$ws.Range("L80").Value = 52997.332
$ws.Range("N80").Value = -54993.332
$ws.Range("H83").Value = 52997.332
$ws.Range("J83").Value = 52997.332
$ws.Range("L83").Value = 158991.996
$ws.Range("N83").Value = -168975.996
$ws.Range("H107").Value = 44228
$ws.Range("J107").Value = 44228
$ws.Range("L107").Value = 44228
$ws.Range("N107").Value = -51908
$ws.Range("H117").Value = 49999.668
$ws.Range("J117").Value = 49999.668
$ws.Range("L117").Value = 49999.668
$ws.Range("N117").Value = -59177.668
$ws.Range("H118").Value = 47498
$ws.Range("J118").Value = 47498
$ws.Range("L118").Value = 47498
$ws.Range("N118").Value = -50812
$ws.Range("H120").Value = 43296
$ws.Range("J120").Value = 43296
$ws.Range("L120").Value = 43296
$ws.Range("N120").Value = -52972
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 46801.668
$ws.Range("J125").Value = 46801.668
$ws.Range("L125").Value = 46801.668
$ws.Range("N125").Value = -56641.668
$ws.Range("H128").Value = 50429
$ws.Range("J128").Value = 50429
$ws.Range("L128").Value = 50429
$ws.Range("N128").Value = -60389
$ws.Range("H131").Value = 50849.5
$ws.Range("J131").Value = 50849.5
$ws.Range("L131").Value = 50849.5
$ws.Range("N131").Value = -60929.5
$ws.Range("H132").Value = 20834886
$ws.Range("I132").Value = 33334318
$ws.Range("K132").Value = 100002954
$ws.Range("M132").Value = -100000424
$ws.Range("H133").Value = 37048.875
$ws.Range("J133").Value = 37048.875
$ws.Range("L133").Value = 37048.875
$ws.Range("N133").Value = -42108.875
$ws.Range("H134").Value = 52285.715
$ws.Range("J134").Value = 52285.715
$ws.Range("L134").Value = 52285.715
$ws.Range("N134").Value = -62425.715
$ws.Range("H135").Value = 44652.637
$ws.Range("J135").Value = 44652.637
$ws.Range("L135").Value = 44652.637
$ws.Range("N135").Value = -54792.637
$ws.Range("H136").Value = 2978.0476
$ws.Range("I136").Value = 1771
$ws.Range("K136").Value = 5313
$ws.Range("M136").Value = -2763

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 433.2
$ws.Range("I94").Value = 385.5625
$ws.Range("K94").Value = 385.5625
$ws.Range("M94").Value = 65.4375
$ws.Range("H117").Value = 48871
$ws.Range("J117").Value = 48871
$ws.Range("L117").Value = 48871
$ws.Range("N117").Value = -58049
$ws.Range("H119").Value = 47257
$ws.Range("J119").Value = 47257
$ws.Range("L119").Value = 47257
$ws.Range("N119").Value = -56933
$ws.Range("H125").Value = 50772
$ws.Range("J125").Value = 50772
$ws.Range("L125").Value = 50772
$ws.Range("N125").Value = -60612
$ws.Range("H134").Value = 3583.15
$ws.Range("I134").Value = 1155.3636
$ws.Range("J134").Value = 4778.9253
$ws.Range("K134").Value = 3466.0908
$ws.Range("L134").Value = 14336.7759
$ws.Range("M134").Value = -931.0907999999999
$ws.Range("N134").Value = -19406.7759
$ws.Range("H141").Value = 31999.818
$ws.Range("J141").Value = 31999.818
$ws.Range("L141").Value = 31999.818
$ws.Range("N141").Value = -42359.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 230542.72
$ws.Range("I31").Value = 2006.8572
$ws.Range("K31").Value = 2006.8572
$ws.Range("M31").Value = -1711.8572
$ws.Range("H34").Value = 230542.72
$ws.Range("I34").Value = 2006.8572
$ws.Range("K34").Value = 2006.8572
$ws.Range("M34").Value = -1804.8572
$ws.Range("H53").Value = 18871
$ws.Range("J53").Value = 18871
$ws.Range("L53").Value = 18871
$ws.Range("N53").Value = -20085
$ws.Range("H99").Value = 1760.1177
$ws.Range("I99").Value = 1340.25
$ws.Range("J99").Value = 2133.3333
$ws.Range("K99").Value = 1340.25
$ws.Range("L99").Value = 2133.3333
$ws.Range("M99").Value = 157.75
$ws.Range("N99").Value = -5129.3333
$ws.Range("H100").Value = 42110
$ws.Range("J100").Value = 42110
$ws.Range("L100").Value = 42110
$ws.Range("N100").Value = -44274
$ws.Range("H118").Value = 44734
$ws.Range("J118").Value = 44734
$ws.Range("L118").Value = 44734
$ws.Range("N118").Value = -48048
$ws.Range("H126").Value = 1760.1177
$ws.Range("I126").Value = 1340.25
$ws.Range("J126").Value = 2133.3333
$ws.Range("K126").Value = 4020.75
$ws.Range("L126").Value = 6399.999899999999
$ws.Range("M126").Value = -1550.75
$ws.Range("N126").Value = -11339.9999
$ws.Range("H132").Value = 52355.57
$ws.Range("I132").Value = 1903.4445
$ws.Range("K132").Value = 5710.333500000001
$ws.Range("M132").Value = -3180.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2517.5881
$ws.Range("I132").Value = 1834.421
$ws.Range("K132").Value = 5503.263
$ws.Range("M132").Value = -2973.263

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 1133.3334
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 1133.3334
$ws.Range("M22").Value = -255
$ws.Range("N22").Value = -1723.3334
$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 550
$ws.Range("J27").Value = 1133.3334
$ws.Range("K27").Value = 550
$ws.Range("L27").Value = 1133.3334
$ws.Range("M27").Value = -443
$ws.Range("N27").Value = -1347.3334
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H68").Value = 2499.8333
$ws.Range("I68").Value = 1914
$ws.Range("J68").Value = 3320
$ws.Range("K68").Value = 1914
$ws.Range("L68").Value = 3320
$ws.Range("M68").Value = -1165
$ws.Range("N68").Value = -4818
$ws.Range("H71").Value = 2499.8333
$ws.Range("I71").Value = 1914
$ws.Range("J71").Value = 3320
$ws.Range("K71").Value = 9570
$ws.Range("L71").Value = 16600
$ws.Range("M71").Value = -5826
$ws.Range("N71").Value = -24088
$ws.Range("H111").Value = 43965.2
$ws.Range("J111").Value = 43965.2
$ws.Range("L111").Value = 43965.2
$ws.Range("N111").Value = -52145.2
$ws.Range("H121").Value = 41137.332
$ws.Range("J121").Value = 41137.332
$ws.Range("L121").Value = 41137.332
$ws.Range("N121").Value = -44631.332
$ws.Range("H122").Value = 2283.3333
$ws.Range("I122").Value = 2263.6365
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6790.9095
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4340.9095
$ws.Range("N122").Value = -12400
$ws.Range("H123").Value = 39413
$ws.Range("J123").Value = 39413
$ws.Range("L123").Value = 39413
$ws.Range("N123").Value = -49213
$ws.Range("H130").Value = 48025
$ws.Range("J130").Value = 48025
$ws.Range("L130").Value = 48025
$ws.Range("N130").Value = -58065

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 44750.5
$ws.Range("J16").Value = 44750.5
$ws.Range("L16").Value = 44750.5
$ws.Range("N16").Value = -45334.5
$ws.Range("H119").Value = 44765.332
$ws.Range("J119").Value = 44765.332
$ws.Range("L119").Value = 44765.332
$ws.Range("N119").Value = -54441.332
$ws.Range("H133").Value = 65614
$ws.Range("J133").Value = 65614
$ws.Range("L133").Value = 65614
$ws.Range("N133").Value = -75734
